# Auto-generated cell updates for Jenova_Profits market-data refresh.
# Each worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) holds Leve profit calcs;
# columns H..N are: currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18
$ws.Cells.Item(18, 8).Value = 1299.5  # H18
$ws.Cells.Item(18, 9).Value = 599  # I18
$ws.Cells.Item(18, 10).Value = 2000  # J18
$ws.Cells.Item(18, 11).Value = 599  # K18
$ws.Cells.Item(18, 12).Value = 2000  # L18
$ws.Cells.Item(18, 13).Value = -315  # M18
$ws.Cells.Item(18, 14).Value = -2568  # N18

# ALC row 118
$ws.Cells.Item(118, 8).Value = 505.55554  # H118
$ws.Cells.Item(118, 9).Value = 505.55554  # I118
$ws.Cells.Item(118, 11).Value = 1516.66662  # K118
$ws.Cells.Item(118, 13).Value = 140.33338  # M118

# ALC row 125
$ws.Cells.Item(125, 8).Value = 3042  # H125
$ws.Cells.Item(125, 9).Value = 946.5  # I125
$ws.Cells.Item(125, 10).Value = 5137.5  # J125
$ws.Cells.Item(125, 11).Value = 8518.5  # K125
$ws.Cells.Item(125, 12).Value = 46237.5  # L125
$ws.Cells.Item(125, 13).Value = -6058.5  # M125
$ws.Cells.Item(125, 14).Value = -51157.5  # N125

# ALC row 132
$ws.Cells.Item(132, 8).Value = 1530.4819  # H132
$ws.Cells.Item(132, 9).Value = 1621.6892  # I132
$ws.Cells.Item(132, 11).Value = 4865.0676  # K132
$ws.Cells.Item(132, 13).Value = -2335.0676  # M132

# ALC row 137
$ws.Cells.Item(137, 8).Value = 5558.067  # H137
$ws.Cells.Item(137, 10).Value = 6928.857  # J137
$ws.Cells.Item(137, 12).Value = 20786.571  # L137
$ws.Cells.Item(137, 14).Value = -25886.571  # N137

# ALC row 138
$ws.Cells.Item(138, 8).Value = 5786.5796  # H138
$ws.Cells.Item(138, 10).Value = 6506.7075  # J138
$ws.Cells.Item(138, 12).Value = 19520.1225  # L138
$ws.Cells.Item(138, 14).Value = -29800.1225  # N138

$ws = $wb.Worksheets.Item("ARM")
# ARM row 26
$ws.Cells.Item(26, 8).Value = 2137.2222  # H26
$ws.Cells.Item(26, 9).Value = 2137.2222  # I26
$ws.Cells.Item(26, 11).Value = 2137.2222  # K26
$ws.Cells.Item(26, 13).Value = -1807.2222  # M26

# ARM row 30
$ws.Cells.Item(30, 8).Value = 1000  # H30
$ws.Cells.Item(30, 9).Value = 1000  # I30
$ws.Cells.Item(30, 11).Value = 1000  # K30
$ws.Cells.Item(30, 13).Value = -850  # M30

# ARM row 45
$ws.Cells.Item(45, 8).Value = 3882.5  # H45
$ws.Cells.Item(45, 9).Value = 2824.1765  # I45
$ws.Cells.Item(45, 10).Value = 5881.5557  # J45
$ws.Cells.Item(45, 11).Value = 2824.1765  # K45
$ws.Cells.Item(45, 12).Value = 5881.5557  # L45
$ws.Cells.Item(45, 13).Value = -2447.1765  # M45
$ws.Cells.Item(45, 14).Value = -6635.5557  # N45

# ARM row 61
$ws.Cells.Item(61, 8).Value = 2714.9832  # H61
$ws.Cells.Item(61, 9).Value = 1663.4255  # I61
$ws.Cells.Item(61, 10).Value = 6833.5835  # J61
$ws.Cells.Item(61, 11).Value = 1663.4255  # K61
$ws.Cells.Item(61, 12).Value = 6833.5835  # L61
$ws.Cells.Item(61, 13).Value = -1451.4255  # M61
$ws.Cells.Item(61, 14).Value = -7257.5835  # N61

# ARM row 96
$ws.Cells.Item(96, 8).Value = 46975  # H96
$ws.Cells.Item(96, 10).Value = 46975  # J96
$ws.Cells.Item(96, 12).Value = 46975  # L96
$ws.Cells.Item(96, 14).Value = -52467  # N96

# ARM row 122
$ws.Cells.Item(122, 8).Value = 4311.625  # H122
$ws.Cells.Item(122, 9).Value = 2727  # I122
$ws.Cells.Item(122, 11).Value = 8181  # K122
$ws.Cells.Item(122, 13).Value = -5731  # M122

# ARM row 132
$ws.Cells.Item(132, 8).Value = 5858.212  # H132
$ws.Cells.Item(132, 9).Value = 4330.647  # I132
$ws.Cells.Item(132, 11).Value = 12991.941  # K132
$ws.Cells.Item(132, 13).Value = -10461.941  # M132

# ARM row 136
$ws.Cells.Item(136, 8).Value = 2714.9832  # H136
$ws.Cells.Item(136, 9).Value = 1663.4255  # I136
$ws.Cells.Item(136, 10).Value = 6833.5835  # J136
$ws.Cells.Item(136, 11).Value = 4990.2765  # K136
$ws.Cells.Item(136, 12).Value = 20500.7505  # L136
$ws.Cells.Item(136, 13).Value = -2440.2765  # M136
$ws.Cells.Item(136, 14).Value = -25600.7505  # N136

# ARM row 139
$ws.Cells.Item(139, 8).Value = 49738.332  # H139
$ws.Cells.Item(139, 10).Value = 49738.332  # J139
$ws.Cells.Item(139, 12).Value = 49738.332  # L139
$ws.Cells.Item(139, 14).Value = -60018.332  # N139

$ws = $wb.Worksheets.Item("BSM")
# BSM row 8
$ws.Cells.Item(8, 8).Value = 2486.25  # H8
$ws.Cells.Item(8, 9).Value = 435  # I8
$ws.Cells.Item(8, 10).Value = 3170  # J8
$ws.Cells.Item(8, 11).Value = 435  # K8
$ws.Cells.Item(8, 12).Value = 3170  # L8
$ws.Cells.Item(8, 13).Value = -295  # M8
$ws.Cells.Item(8, 14).Value = -3450  # N8

# BSM row 17
$ws.Cells.Item(17, 8).Value = 860  # H17
$ws.Cells.Item(17, 10).Value = 860  # J17
$ws.Cells.Item(17, 12).Value = 860  # L17
$ws.Cells.Item(17, 14).Value = -1204  # N17

# BSM row 29
$ws.Cells.Item(29, 8).Value = 499  # H29
$ws.Cells.Item(29, 9).Value = 499  # I29
$ws.Cells.Item(29, 10).Value = 0  # J29
$ws.Cells.Item(29, 11).Value = 499  # K29
$ws.Cells.Item(29, 13).Value = -210  # M29
$ws.Cells.Item(29, 14).ClearContents()  # N29 removed

# BSM row 35
$ws.Cells.Item(35, 8).Value = 103979.8  # H35
$ws.Cells.Item(35, 10).Value = 107474.75  # J35
$ws.Cells.Item(35, 12).Value = 107474.75  # L35
$ws.Cells.Item(35, 14).Value = -108094.75  # N35

# BSM row 99
$ws.Cells.Item(99, 8).Value = 5315.1055  # H99
$ws.Cells.Item(99, 10).Value = 4798.3335  # J99
$ws.Cells.Item(99, 12).Value = 4798.3335  # L99
$ws.Cells.Item(99, 14).Value = -7794.3335  # N99

# BSM row 134
$ws.Cells.Item(134, 8).Value = 17628.928  # H134
$ws.Cells.Item(134, 9).Value = 2028.4117  # I134
$ws.Cells.Item(134, 10).Value = 61830.39  # J134
$ws.Cells.Item(134, 11).Value = 6085.2351  # K134
$ws.Cells.Item(134, 12).Value = 185491.17  # L134
$ws.Cells.Item(134, 13).Value = -3550.2351  # M134
$ws.Cells.Item(134, 14).Value = -190561.17  # N134

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Cells.Item(22, 8).Value = 544.26666  # H22
$ws.Cells.Item(22, 9).Value = 317  # I22
$ws.Cells.Item(22, 11).Value = 317  # K22
$ws.Cells.Item(22, 13).Value = 33  # M22

# CRP row 41
$ws.Cells.Item(41, 8).Value = 35740  # H41
$ws.Cells.Item(41, 10).Value = 54866.668  # J41
$ws.Cells.Item(41, 12).Value = 54866.668  # L41
$ws.Cells.Item(41, 14).Value = -55722.668  # N41

# CRP row 110
$ws.Cells.Item(110, 8).Value = 0  # H110
$ws.Cells.Item(110, 10).Value = 0  # J110
$ws.Cells.Item(110, 14).ClearContents()  # N110 removed

# CRP row 141
$ws.Cells.Item(141, 8).Value = 100427.86  # H141
$ws.Cells.Item(141, 9).Value = 25000  # I141
$ws.Cells.Item(141, 11).Value = 25000  # K141
$ws.Cells.Item(141, 13).Value = -19820  # M141

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Cells.Item(5, 8).Value = 91166.336  # H5
$ws.Cells.Item(5, 10).Value = 3166.1667  # J5
$ws.Cells.Item(5, 12).Value = 9498.500100000001  # L5
$ws.Cells.Item(5, 14).Value = -9722.500100000001  # N5

# CUL row 17
$ws.Cells.Item(17, 8).Value = 2463.2  # H17
$ws.Cells.Item(17, 9).Value = 3888  # I17
$ws.Cells.Item(17, 11).Value = 11664  # K17
$ws.Cells.Item(17, 13).Value = -11495  # M17

# CUL row 56
$ws.Cells.Item(56, 8).Value = 7239.5  # H56
$ws.Cells.Item(56, 9).Value = 7239.5  # I56
$ws.Cells.Item(56, 11).Value = 7239.5  # K56
$ws.Cells.Item(56, 13).Value = -6709.5  # M56

# CUL row 86
$ws.Cells.Item(86, 8).Value = 294.5  # H86
$ws.Cells.Item(86, 10).Value = 294.5  # J86
$ws.Cells.Item(86, 12).Value = 883.5  # L86
$ws.Cells.Item(86, 14).Value = -3255.5  # N86

# CUL row 89
$ws.Cells.Item(89, 8).Value = 294.5  # H89
$ws.Cells.Item(89, 10).Value = 294.5  # J89
$ws.Cells.Item(89, 12).Value = 2650.5  # L89
$ws.Cells.Item(89, 14).Value = -14506.5  # N89

# CUL row 135
$ws.Cells.Item(135, 8).Value = 91166.336  # H135
$ws.Cells.Item(135, 10).Value = 3166.1667  # J135
$ws.Cells.Item(135, 12).Value = 28495.5003  # L135
$ws.Cells.Item(135, 14).Value = -33565.5003  # N135

# CUL row 136
$ws.Cells.Item(136, 8).Value = 7909.8  # H136
$ws.Cells.Item(136, 9).Value = 6887.25  # I136
$ws.Cells.Item(136, 10).Value = 12000  # J136
$ws.Cells.Item(136, 11).Value = 20661.75  # K136
$ws.Cells.Item(136, 12).Value = 36000  # L136
$ws.Cells.Item(136, 13).Value = -15561.75  # M136
$ws.Cells.Item(136, 14).Value = -46200  # N136

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Cells.Item(122, 8).Value = 3830.4211  # H122
$ws.Cells.Item(122, 9).Value = 3230.2727  # I122
$ws.Cells.Item(122, 10).Value = 4655.625  # J122
$ws.Cells.Item(122, 11).Value = 9690.8181  # K122
$ws.Cells.Item(122, 12).Value = 13966.875  # L122
$ws.Cells.Item(122, 13).Value = -7240.8181  # M122
$ws.Cells.Item(122, 14).Value = -18866.875  # N122

# GSM row 132
$ws.Cells.Item(132, 8).Value = 189283.31  # H132
$ws.Cells.Item(132, 9).Value = 229722.8  # I132
$ws.Cells.Item(132, 10).Value = 70660.8  # J132
$ws.Cells.Item(132, 11).Value = 689168.3999999999  # K132
$ws.Cells.Item(132, 12).Value = 211982.4  # L132
$ws.Cells.Item(132, 13).Value = -686638.3999999999  # M132
$ws.Cells.Item(132, 14).Value = -217042.4  # N132

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Cells.Item(7, 8).Value = 377473.97  # H7
$ws.Cells.Item(7, 9).Value = 7581.357  # I7
$ws.Cells.Item(7, 10).Value = 775819.9  # J7
$ws.Cells.Item(7, 11).Value = 7581.357  # K7
$ws.Cells.Item(7, 12).Value = 775819.9  # L7
$ws.Cells.Item(7, 13).Value = -7469.357  # M7
$ws.Cells.Item(7, 14).Value = -776043.9  # N7

# LTW row 55
$ws.Cells.Item(55, 8).Value = 1548.2354  # H55
$ws.Cells.Item(55, 9).Value = 700.4167  # I55
$ws.Cells.Item(55, 10).Value = 3583  # J55
$ws.Cells.Item(55, 11).Value = 700.4167  # K55
$ws.Cells.Item(55, 12).Value = 3583  # L55
$ws.Cells.Item(55, 13).Value = -527.4167  # M55
$ws.Cells.Item(55, 14).Value = -3929  # N55

# LTW row 61
$ws.Cells.Item(61, 8).Value = 3593  # H61
$ws.Cells.Item(61, 9).Value = 1916.45  # I61
$ws.Cells.Item(61, 10).Value = 5269.55  # J61
$ws.Cells.Item(61, 11).Value = 1916.45  # K61
$ws.Cells.Item(61, 12).Value = 5269.55  # L61
$ws.Cells.Item(61, 13).Value = -1714.45  # M61
$ws.Cells.Item(61, 14).Value = -5673.55  # N61

# LTW row 113
$ws.Cells.Item(113, 8).Value = 3593  # H113
$ws.Cells.Item(113, 9).Value = 1916.45  # I113
$ws.Cells.Item(113, 10).Value = 5269.55  # J113
$ws.Cells.Item(113, 11).Value = 1916.45  # K113
$ws.Cells.Item(113, 12).Value = 5269.55  # L113
$ws.Cells.Item(113, 13).Value = 253.55  # M113
$ws.Cells.Item(113, 14).Value = -9609.549999999999  # N113

# LTW row 126
$ws.Cells.Item(126, 8).Value = 377473.97  # H126
$ws.Cells.Item(126, 9).Value = 7581.357  # I126
$ws.Cells.Item(126, 10).Value = 775819.9  # J126
$ws.Cells.Item(126, 11).Value = 22744.071  # K126
$ws.Cells.Item(126, 12).Value = 2327459.7  # L126
$ws.Cells.Item(126, 13).Value = -20274.071  # M126
$ws.Cells.Item(126, 14).Value = -2332399.7  # N126

# LTW row 132
$ws.Cells.Item(132, 8).Value = 2951.549  # H132
$ws.Cells.Item(132, 9).Value = 2481.5715  # I132
$ws.Cells.Item(132, 11).Value = 7444.7145  # K132
$ws.Cells.Item(132, 13).Value = -4914.7145  # M132

# LTW row 136
$ws.Cells.Item(136, 8).Value = 3389  # H136
$ws.Cells.Item(136, 9).Value = 2738.1765  # I136
$ws.Cells.Item(136, 10).Value = 4311  # J136
$ws.Cells.Item(136, 11).Value = 8214.529500000001  # K136
$ws.Cells.Item(136, 12).Value = 12933  # L136
$ws.Cells.Item(136, 13).Value = -5664.529500000001  # M136
$ws.Cells.Item(136, 14).Value = -18033  # N136

$ws = $wb.Worksheets.Item("WVR")
# WVR row 40
$ws.Cells.Item(40, 8).Value = 94166.664  # H40
$ws.Cells.Item(40, 10).Value = 86000  # J40
$ws.Cells.Item(40, 12).Value = 86000  # L40
$ws.Cells.Item(40, 14).Value = -86298  # N40

# WVR row 132
$ws.Cells.Item(132, 8).Value = 21188.018  # H132
$ws.Cells.Item(132, 9).Value = 1840.0256  # I132
$ws.Cells.Item(132, 10).Value = 71492.8  # J132
$ws.Cells.Item(132, 11).Value = 5520.0768  # K132
$ws.Cells.Item(132, 12).Value = 214478.4  # L132
$ws.Cells.Item(132, 13).Value = -2990.0768  # M132
$ws.Cells.Item(132, 14).Value = -219538.4  # N132

# WVR row 133
$ws.Cells.Item(133, 8).Value = 78997.5  # H133
$ws.Cells.Item(133, 10).Value = 78997.5  # J133
$ws.Cells.Item(133, 12).Value = 78997.5  # L133
$ws.Cells.Item(133, 14).Value = -89117.5  # N133

# WVR row 136
$ws.Cells.Item(136, 8).Value = 56347.727  # H136
$ws.Cells.Item(136, 9).Value = 11300  # I136
$ws.Cells.Item(136, 10).Value = 672000  # J136
$ws.Cells.Item(136, 11).Value = 33900  # K136
$ws.Cells.Item(136, 12).Value = 2016000  # L136
$ws.Cells.Item(136, 13).Value = -31350  # M136
$ws.Cells.Item(136, 14).Value = -2021100  # N136
